# Update Fonds de solidarite Volet 2 regional data with 2020-07-22 figures.
# For each affected row, nombre_aides (column C) and montant_total (column D)
# are refreshed with the latest counts/amounts. Values are written with a
# leading apostrophe so Excel stores them as text (matching the source
# data's text-typed numeric columns) instead of auto-converting to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = "'131"
$ws.Range("D9").Value = "'298238.55"
$ws.Range("C10").Value = "'281"
$ws.Range("D10").Value = "'867021.77"
$ws.Range("C11").Value = "'120"
$ws.Range("D11").Value = "'341125.40"
$ws.Range("C15").Value = "'419"
$ws.Range("D15").Value = "'1364254.00"
$ws.Range("C19").Value = "'143"
$ws.Range("D19").Value = "'346569.00"
$ws.Range("C20").Value = "'291"
$ws.Range("D20").Value = "'880952.97"
$ws.Range("C27").Value = "'187"
$ws.Range("D27").Value = "'414926.00"
$ws.Range("C29").Value = "'387"
$ws.Range("D29").Value = "'1181532.00"
$ws.Range("C31").Value = "'301"
$ws.Range("D31").Value = "'825218.74"
$ws.Range("C38").Value = "'47"
$ws.Range("D38").Value = "'117439.00"
$ws.Range("C39").Value = "'41"
$ws.Range("D39").Value = "'192229.92"
$ws.Range("C40").Value = "'68"
$ws.Range("D40").Value = "'301198.15"
$ws.Range("C41").Value = "'4"
$ws.Range("D41").Value = "'21132.00"
$ws.Range("C42").Value = "'223"
$ws.Range("D42").Value = "'584420.74"
$ws.Range("C44").Value = "'454"
$ws.Range("D44").Value = "'1506039.81"
$ws.Range("C45").Value = "'314"
$ws.Range("D45").Value = "'969653.79"
$ws.Range("C47").Value = "'22"
$ws.Range("D47").Value = "'110220.65"
$ws.Range("C52").Value = "'2847"
$ws.Range("D52").Value = "'6866016.41"
$ws.Range("C78").Value = "'175"
$ws.Range("D78").Value = "'399736.00"
$ws.Range("C80").Value = "'437"
$ws.Range("D80").Value = "'1318734.92"
$ws.Range("C81").Value = "'167"
$ws.Range("D81").Value = "'471749.09"
$ws.Range("C84").Value = "'409"
$ws.Range("D84").Value = "'886650.00"
$ws.Range("C86").Value = "'909"
$ws.Range("D86").Value = "'2598207.86"
$ws.Range("C88").Value = "'838"
$ws.Range("D88").Value = "'2216502.50"
$ws.Range("C90").Value = "'33"
$ws.Range("D90").Value = "'85000.00"
